# Update the "Adam / 12288 / seed 0" row (row 8) of the seed-investigation
# report with the re-measured accuracy numbers, letting the Diff. column
# formula (G8 = E8-F8) recalculate automatically, then leave the selection
# on F9 as it was when the report was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 0.502
$ws.Range("F8").Value = 0.543

$ws.Range("F9").Select()
